# Added two new Mac-Addresses: append 10 new device rows (147-156) to the
# master-reg_center_device test-data sheet, following the existing pattern
# (regcntr_id constant, device_id incrementing, lang/active/cr_by/cr_dtimes
# repeating the same values used throughout the table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147
$startDeviceId = 3000166
$rowCount = 10

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $startDeviceId + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Match the saved selection after entering the new data.
$ws.Range("C152").Select()
